$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (2023-10-13, Provincia de Linares) was
# recorded ahead of the existing rows; insert a row at 74 so the later
# entries (old 74-78) shift down to 75-79, then populate the new row.
$ws.Rows("74:74").Insert()

$ws.Range("A74").Value = 3
$ws.Range("B74").Value = "Femacal de La Calera"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = 45212
$ws.Range("E74").Value = 5
$ws.Range("F74").Value = 300000000
$ws.Range("G74").Value = "Espárragos"
$ws.Range("H74").Value = "Verde"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 1800
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = 1800
$ws.Range("N74").Value = "$/kilo"
$ws.Range("O74").Value = "Provincia de Linares"
$ws.Range("P74").Value = 1800
$ws.Range("Q74").Value = 1
$ws.Range("R74").Value = "Hortaliza"
